$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F35").Value = -2.474451917775058
$ws.Range("F36").Value = -2.4814437950556
$ws.Range("F37").Value = -2.489118030193063
$ws.Range("F38").Value = -2.497474623187442
$ws.Range("F39").Value = -2.504901719068524
$ws.Range("F40").Value = -2.512140735605447
$ws.Range("F42").Value = -2.528606782337825
$ws.Range("F43").Value = -2.535598659618369
$ws.Range("F44").Value = -2.54352003401221
$ws.Range("F45").Value = -2.551629487750215
$ws.Range("F46").Value = -2.559056583631296
$ws.Range("F47").Value = -2.566542739424595
$ws.Range("F48").Value = -2.574652193162594
$ws.Range("F49").Value = -2.582761646900595
$ws.Range("F51").Value = -2.597862977919144
$ws.Range("F52").Value = -2.605784352312977
$ws.Range("F53").Value = -2.613893806050984
$ws.Range("F54").Value = -2.621320901932065
$ws.Range("F81").Value = -2.558505259722823
$ws.Range("F82").Value = -2.613125735955284
$ws.Range("F83").Value = -2.670506987783486
$ws.Range("F85").Value = -2.768461063707048
$ws.Range("F86").Value = -2.812082632340831
$ws.Range("F87").Value = -2.856930867366638
$ws.Range("F92").Value = -2.645096737470864
$ws.Range("F93").Value = -2.699096597998986
$ws.Range("F94").Value = -2.754692000339417
$ws.Range("F95").Value = -2.812947599279458
$ws.Range("F96").Value = -2.868480099813693
$ws.Range("F97").Value = -2.913253757132945
$ws.Range("F98").Value = -2.955945917568026
$ws.Range("F99").Value = -3.00094214162047
$ws.Range("F103").Value = -2.713915989792761
$ws.Range("F104").Value = -2.767839347684372
$ws.Range("F105").Value = -2.822615060882911
$ws.Range("F106").Value = -2.878658774433729
$ws.Range("F107").Value = -2.935421146364162
$ws.Range("F108").Value = -2.986578757916763
$ws.Range("F109").Value = -3.028950551513024
$ws.Range("F110").Value = -3.072742249194393
$ws.Range("F111").Value = -3.116872129177024
$ws.Range("F114").Value = -2.754578114113495
$ws.Range("F115").Value = -2.807845289465562
$ws.Range("F116").Value = -2.859121085979382
$ws.Range("F117").Value = -2.917666950668647
$ws.Range("F118").Value = -2.973610553659908
$ws.Range("F119").Value = -3.031720456456978
$ws.Range("F120").Value = -3.076424275331069
$ws.Range("F121").Value = -3.119548000777475
$ws.Range("F122").Value = -3.163582216167744
$ws.Range("F125").Value = -2.788584805932618
$ws.Range("F126").Value = -2.841858066267444
$ws.Range("F127").Value = -2.896101910275515
$ws.Range("F128").Value = -2.950707437300936
$ws.Range("F129").Value = -3.006470958170002
$ws.Range("F130").Value = -3.062727936360075
$ws.Range("F131").Value = -3.114135288123573
$ws.Range("F132").Value = -3.157012822852406
$ws.Range("F133").Value = -3.20053946384706
$ws.Range("F136").Value = -2.812227221411305
$ws.Range("F137").Value = -2.864706399126527
$ws.Range("F138").Value = -2.916983844674297
$ws.Range("F139").Value = -2.972565592774685
$ws.Range("F140").Value = -3.027487821381273
$ws.Range("F141").Value = -3.084367390289853
$ws.Range("F142").Value = -3.141838906964608
$ws.Range("F143").Value = -3.184269154705247
$ws.Range("F144").Value = -3.227776674556784
$ws.Range("F147").Value = -2.832489980898743
$ws.Range("F148").Value = -2.883972078524685
$ws.Range("F149").Value = -2.936730587200261
$ws.Range("F150").Value = -2.988663485158066
$ws.Range("F151").Value = -3.046076118672538
$ws.Range("F152").Value = -3.102722513496122
$ws.Range("F153").Value = -3.161612226370578
$ws.Range("F154").Value = -3.207227486688069
$ws.Range("F155").Value = -3.25030299440986
$ws.Range("F158").Value = -2.847618439783126
$ws.Range("F159").Value = -2.899033232366425
$ws.Range("F160").Value = -2.950200078939421
$ws.Range("F161").Value = -3.003212002989778
$ws.Range("F162").Value = -3.058908103099831
$ws.Range("F163").Value = -3.115855674263634
$ws.Range("F164").Value = -3.174906054089296
$ws.Range("F165").Value = -3.225024738811097
$ws.Range("F166").Value = -3.267659759314567
$ws.Range("F169").Value = -2.861082916410609
$ws.Range("F170").Value = -2.912124220513491
$ws.Range("F171").Value = -2.96218244988479
$ws.Range("F172").Value = -3.01557492310028
$ws.Range("F173").Value = -3.070320051149699
$ws.Range("F174").Value = -3.127437491041535
$ws.Range("F175").Value = -3.186620212211471
$ws.Range("F176").Value = -3.240642990738449
$ws.Range("F177").Value = -3.283257311464205
$ws.Range("F180").Value = -2.871333072897983
$ws.Range("F181").Value = -2.920218306848614
$ws.Range("F182").Value = -2.972251535708022
$ws.Range("F183").Value = -3.024248683206123
$ws.Range("F184").Value = -3.078909829184266
$ws.Range("F185").Value = -3.137040609887663
$ws.Range("F186").Value = -3.196715147652993
$ws.Range("F187").Value = -3.253306900269987
$ws.Range("F188").Value = -3.295585249564823
$ws.Range("F191").Value = -2.880726833488755
$ws.Range("F192").Value = -2.929079179246594
$ws.Range("F193").Value = -2.979511554689026
$ws.Range("F194").Value = -3.032701747877713
$ws.Range("F195").Value = -3.087793898648273
$ws.Range("F196").Value = -3.144340303716703
$ws.Range("F197").Value = -3.203421241999518
$ws.Range("F198").Value = -3.264957319653926
$ws.Range("F199").Value = -3.306662230036727
$ws.Range("F202").Value = -2.887976249321166
$ws.Range("F203").Value = -2.936444977148717
$ws.Range("F204").Value = -2.98618627581132
$ws.Range("F205").Value = -3.038414680868706
$ws.Range("F206").Value = -3.093464661611373
$ws.Range("F207").Value = -3.151109305509923
$ws.Range("F208").Value = -3.210782470897806
$ws.Range("F209").Value = -3.272718795372747
$ws.Range("F210").Value = -3.315805959360634
Write-Host "Applied 122 changes to column F"
